$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price, "Volume(1h)") target values for the refreshed cryptos list.
# D-column values that look like plain decimal numbers are given a leading
# apostrophe (exactly like a user typing '0.9998 into Excel) so the cell keeps
# storing them as text instead of Excel auto-converting them to numbers.
$data = @{
    2 = @("30.275.85", "  +0.69%  ")
    3 = @("1.865.03", "  +0.33%  ")
    4 = @("'0.9998", "  -0.07%  ")
    5 = @("'237.23", "  +1.83%  ")
    6 = @("'0.9996", "  -0.08%  ")
    7 = @("'0.4682", "  +0.63%  ")
    8 = @("'0.2860", "  +2.18%  ")
    9 = @("'0.06542", "  +0.29%  ")
    10 = @("'22.35", "  +14.89%  ")
    11 = @("'0.07912", "  +1.32%  ")
    12 = @("'97.75", "  +1.61%  ")
    13 = @("1.869.44", "  +0.60%  ")
    14 = @("'5.176", "  +1.44%  ")
    15 = @("'0.6844", "  +3.29%  ")
    16 = @("'279.61", "  -0.07%  ")
    17 = @("30.266.50", "  +0.55%  ")
    18 = @("'13.69", "  +9.25%  ")
    19 = @("'0.9997", "  -0.08%  ")
    20 = @("'5.397", "  -1.32%  ")
    21 = @("'0.000007334", "  +1.83%  ")
    22 = @("2.112.07", "  -0.03%  ")
    23 = @("'1.000", "  -0.03%  ")
    24 = @("'6.176", "  +1.24%  ")
    25 = @("'168.16", "  +0.79%  ")
    26 = @("'9.266", "  -0.20%  ")
    27 = @("'19.15", "  +1.99%  ")
    28 = @("'1.938", "  +2.10%  ")
    29 = @("'1.381", "  +3.95%  ")
    30 = @("'0.09820", "  +3.16%  ")
    31 = @("'4.402", "  -0.21%  ")
    32 = @("'1.480", "  +1.28%  ")
    33 = @("'4.070", "  -0.07%  ")
    34 = @("'0.04750", "  +2.94%  ")
    35 = @("'1.138", "  +4.37%  ")
    36 = @("'0.7115", "  +1.93%  ")
    37 = @("'2.704", "  +0.29%  ")
    38 = @("'0.01875", "  +1.96%  ")
    39 = @("'2.615", "  +4.42%  ")
    40 = @("'76.65", "  +6.10%  ")
    41 = @("'6.313", "  +0.75%  ")
    42 = @("'1.959", "  +2.94%  ")
    43 = @("'0.8509", "  -0.17%  ")
    44 = @("'0.4191", "  +1.43%  ")
    45 = @("'0.9991", "  -0.12%  ")
    46 = @("'103.51", "  -0.36%  ")
    47 = @("'968.50", "  -2.81%  ")
    48 = @("'7.240", "  +1.15%  ")
    49 = @("'9.319", "  +1.04%  ")
    50 = @("'34.23", "  +1.03%  ")
    51 = @("'0.05645", "  +0.44%  ")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("D$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
}
